$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense data, most-recent-first, replacing/extending the old
# dinner/lunch rows (rows 2-3) with the fuller history (rows 2-16).
$data = @(
    @("Dal bati", 10,      45825.54185799768),
    @("fghj",     55,      45825.00011574074),
    @("tttt",     4444444, 45825.00011574074),
    @("dfghj",    122222,  45825.00011574074),
    @("Dal bati", 10000,   45824.54185799768),
    @("poiuyt",   7,       45824.00011574074),
    @("dinner",   1000,    45817.54185799768),
    @("lunch",    500,     45817.54185799768),
    @("lunch",    500,     45817.54185799768),
    @("test",     500,     45816.54185799768),
    @("test",     5000,    45814.54185799768),
    @("djhbf",    5000,    45813.54185799768),
    @("dfghjk",   567777,  45811.00011574074),
    @("djhbf",    5000,    45810.54185799768),
    @("toffee",   100,     45809.54185799768)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}

# The existing date cell C2 already carries the mm/dd/yyyy-style date
# formatting (style index 1 / numFmtId 14). Propagate that same format
# to all the newly added date cells instead of creating a new number
# format entry.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2:C16").PasteSpecial(-4122) | Out-Null
